# Auto-generated edit script: updates crypto price/volume table cells
# to match the target snapshot (GitHub Actions crypto-list refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D sometimes hold digit-and-dot strings (e.g. "240.30", "1.00")
# that Excel would otherwise auto-coerce into numbers (dropping trailing
# zeros / renormalizing). Force those specific cells to Text format first
# so the literal string is preserved exactly, matching the source data.

$textCells = @('D5', 'D6', 'D12', 'D14', 'D18', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D26', 'D27', 'D28', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '95.834.68'
$ws.Range('E2').Value = '  +3.57%  '
$ws.Range('D3').Value = '3.613.91'
$ws.Range('E3').Value = '  +6.56%  '
$ws.Range('D5').Value = '240.30'
$ws.Range('E5').Value = '  +4.78%  '
$ws.Range('D6').Value = '654.16'
$ws.Range('E6').Value = '  +6.20%  '
$ws.Range('E7').Value = '  +7.87%  '
$ws.Range('E8').Value = '  +5.55%  '
$ws.Range('E10').Value = '  +6.15%  '
$ws.Range('D11').Value = '3.614.30'
$ws.Range('E11').Value = '  +6.62%  '
$ws.Range('D12').Value = '43.43'
$ws.Range('E12').Value = '  +2.10%  '
$ws.Range('E13').Value = '  +2.28%  '
$ws.Range('D14').Value = '6.37'
$ws.Range('E14').Value = '  +2.78%  '
$ws.Range('D15').Value = '4.302.09'
$ws.Range('E15').Value = '  +6.36%  '
$ws.Range('D16').Value = '95.746.35'
$ws.Range('E16').Value = '  +3.57%  '
$ws.Range('E17').Value = '  +5.56%  '
$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').Value = '8.25'
$ws.Range('E18').Value = '  +3.04%  '
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').Value = '3.618.06'
$ws.Range('E19').Value = '  +6.81%  '
$ws.Range('D20').Value = '12.59'
$ws.Range('E20').Value = '  +10.31%  '
$ws.Range('D21').Value = '18.20'
$ws.Range('E21').Value = '  +2.96%  '
$ws.Range('D22').Value = '3.54'
$ws.Range('E22').Value = '  +7.49%  '
$ws.Range('D23').Value = '0.488'
$ws.Range('E23').Value = '  +13.34%  '
$ws.Range('D24').Value = '512.87'
$ws.Range('E24').Value = '  +3.99%  '
$ws.Range('D25').Value = '0.0000198'
$ws.Range('E25').Value = '  +8.64%  '
$ws.Range('D26').Value = '6.70'
$ws.Range('E26').Value = '  +3.38%  '
$ws.Range('D27').Value = '97.32'
$ws.Range('E27').Value = '  +8.02%  '
$ws.Range('D28').Value = '12.88'
$ws.Range('E28').Value = '  +8.58%  '
$ws.Range('B29').Value = 'WrappedeETH'
$ws.Range('C29').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D29').Value = '3.799.12'
$ws.Range('E29').Value = '  +6.23%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value = '3.22'
$ws.Range('E30').Value = '  +19.64%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').Value = '11.39'
$ws.Range('E31').Value = '  +2.12%  '
$ws.Range('B32').Value = 'Dai'
$ws.Range('C32').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  -0.04%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '0.140'
$ws.Range('E33').Value = '  +5.00%  '
$ws.Range('B34').Value = 'Binance-PegBSC-USD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D34').Value = '0.993'
$ws.Range('E34').Value = '  +0.51%  '
$ws.Range('B35').Value = 'Cronos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D35').Value = '0.179'
$ws.Range('E35').Value = '  +4.40%  '
$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D36').Value = '32.02'
$ws.Range('E36').Value = '  +9.05%  '
$ws.Range('B37').Value = 'PolygonEcosystemToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D37').Value = '0.566'
$ws.Range('E37').Value = '  +5.64%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D38').Value = '8.42'
$ws.Range('E38').Value = '  +13.72%  '
$ws.Range('B39').Value = 'Bittensor'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D39').Value = '570.90'
$ws.Range('E39').Value = '  +3.75%  '
$ws.Range('B40').Value = 'Fetch.AI'
$ws.Range('C40').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D40').Value = '1.52'
$ws.Range('E40').Value = '  +10.19%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '0.152'
$ws.Range('E41').Value = '  +2.22%  '
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').Value = '0.936'
$ws.Range('E42').Value = '  +3.49%  '
$ws.Range('B43').Value = 'USDe'
$ws.Range('C43').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  -0.01%  '
$ws.Range('B44').Value = 'ImmutableX'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D44').Value = '1.75'
$ws.Range('E44').Value = '  +2.65%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '34.65'
$ws.Range('E45').Value = '  +40.94%  '
$ws.Range('B46').Value = 'Filecoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D46').Value = '5.76'
$ws.Range('E46').Value = '  +6.48%  '
$ws.Range('B47').Value = 'WhiteBITCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D47').Value = '23.79'
$ws.Range('E47').Value = '  +0.73%  '
$ws.Range('D48').Value = '0.0421'
$ws.Range('E48').Value = '  +4.92%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').Value = '2.26'
$ws.Range('E49').Value = '  +8.40%  '
$ws.Range('B50').Value = 'OKB'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D50').Value = '54.38'
$ws.Range('E50').Value = '  +2.28%  '
$ws.Range('D51').Value = '3.47'
$ws.Range('E51').Value = '  -5.42%  '
